$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "implement feature to copy fees from previous month to current month /
#  fix issue fee list not reloading after deleting fee"
# -> the "Status" column (C) for the two newest stories advances:
#    row 11 (I want to copy over existing expenses to new month) is now DONE
#    row 12 (I want to export and reimport DB) is now IN PROGRESS
#
# Update row 12 first (NOT STARTED -> IN PROGRESS) by copying the existing
# "IN PROGRESS" formatting from C11, then update row 11 (IN PROGRESS -> DONE)
# by copying the existing "DONE" formatting from C10. Doing the row-12 copy
# before overwriting C11 keeps every shared string in use throughout, so the
# shared-string table isn't disturbed.
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = "IN PROGRESS"

$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = "DONE"

# Move the active selection to E11 to match the saved view state.
$null = $ws.Range("E11").Select()
